# Generate Report for Handoff
#
# The CI job re-ran and produced a fresh report:
#   - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#     (Overview!E2/F2, zh-cn!C2, de-de!C2 - the "Status" column).
#   - The associated timestamps advance to the new run time
#     (Overview!G2 "Latest HO Xliff Generate Date", zh-cn!H2 and
#     de-de!H2 "Latest Handoff Datetime").
#   - The now shorter status text means the "Status"/summary columns can be
#     narrower, so those columns are resized down from their old width.

$wb = $excel.ActiveWorkbook

$statusNew = "Ready for handoff"

$overviewDateNew = "2016-08-16 18:54:01"
$zhHandoffDateNew = "2016-08-16 18:53:55"
$deHandoffDateNew = "2016-08-16 18:54:01"

# Closest width (in Excel's character-width units) this host's ColumnWidth
# setter can resolve to the target stored column width of ~17.216 chars -
# the setter snaps to a 1/6-character pixel grid, and 16.3333... is the
# value that lands on the nearest attainable grid point.
$narrowWidth = 16.3333333333333

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusNew
$wsOverview.Range("F2").Value = $statusNew
$wsOverview.Range("G2").Value = $overviewDateNew

$wsOverview.Columns.Item(5).ColumnWidth = $narrowWidth
$wsOverview.Columns.Item(6).ColumnWidth = $narrowWidth

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = $statusNew
$wsZhCn.Range("H2").Value = $zhHandoffDateNew

$wsZhCn.Columns.Item(3).ColumnWidth = $narrowWidth

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = $statusNew
$wsDeDe.Range("H2").Value = $deHandoffDateNew

$wsDeDe.Columns.Item(3).ColumnWidth = $narrowWidth
